$wb = $excel.ActiveWorkbook

# Update the neighborhood labels on the main "db_wijkconfig" sheet:
# remove the colon after "Buurt A"/"Buurt B" (e.g. "Buurt A: reguliere buurt" -> "Buurt A reguliere buurt")
$ws1 = $wb.Worksheets.Item("db_wijkconfig")
$ws1.Range("B2").Value = "Buurt A reguliere buurt"
$ws1.Range("B3").Value = "Buurt B warmtenet"

# Reflect the active selection left on this sheet when the file was saved
$ws1.Range("E10").Select()
